$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 297
$ws1.Range("F6").Value = 388
$ws1.Range("F7").Value = 845
$ws1.Range("F9").Value = 499
$ws1.Range("F11").Value = 290
$ws1.Range("F12").Value = 124
$ws1.Range("F13").Value = 104
$ws1.Range("F15").Value = 26
$ws1.Range("F16").Value = 404
$ws1.Range("F17").Value = 6567
$ws1.Range("F20").Value = 19
$ws1.Range("F21").Value = 7509
$ws1.Range("F24").Value = 3373
$ws1.Range("F25").Value = 21
$ws1.Range("F26").Value = 1161
$ws1.Range("F27").Value = 877
$ws1.Range("F29").Value = 18
$ws1.Range("F31").Value = 65
$ws1.Range("F33").Value = 190
$ws1.Range("F34").Value = 1593
$ws1.Range("F37").Value = 51
$ws1.Range("F39").Value = 1176
$ws1.Range("F40").Value = 1692
$ws1.Range("F41").Value = 2126

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 47

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1214

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1214
$ws4.Range("F7").Value = 297
$ws4.Range("F8").Value = 388
$ws4.Range("F9").Value = 845
$ws4.Range("F11").Value = 499
$ws4.Range("F14").Value = 290
$ws4.Range("F15").Value = 124
$ws4.Range("F17").Value = 104
$ws4.Range("F19").Value = 26
$ws4.Range("F20").Value = 404
$ws4.Range("F21").Value = 6567
$ws4.Range("F24").Value = 19
$ws4.Range("F25").Value = 7509
$ws4.Range("F28").Value = 3373
$ws4.Range("F29").Value = 21
$ws4.Range("F30").Value = 1161
$ws4.Range("F31").Value = 877
$ws4.Range("F33").Value = 18
$ws4.Range("F35").Value = 65
$ws4.Range("F36").Value = 47
$ws4.Range("F38").Value = 190
$ws4.Range("F39").Value = 1594
$ws4.Range("F42").Value = 51
$ws4.Range("F44").Value = 1176
$ws4.Range("F45").Value = 1692
$ws4.Range("F47").Value = 2126
